$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    # Route the literal string through a formula + paste-as-values round
    # trip so Excel's "looks like a date" auto-conversion never kicks in
    # (a plain .Value assignment would silently turn e.g. "11/29/2024"
    # into a date serial number).
    $cell = $ws.Cells.Item($row, $col)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163) | Out-Null
}

# Row 15
Set-TextValue 15 1 "4CH3Z"
Set-TextValue 15 2 "11/29/2024"
Set-TextValue 15 3 "Ashar Nadeem"
Set-TextValue 15 4 "0322-7287568"
Set-TextValue 15 5 "Lahore"
Set-TextValue 15 6 "0322-7287568"
Set-TextValue 15 7 "xxdxsdxdsx"
Set-TextValue 15 8 "csssdccs"
Set-TextValue 15 9 "cscsd"

# Row 16
Set-TextValue 16 1 "BXMIY"
Set-TextValue 16 2 "11/29/2024"
Set-TextValue 16 3 "Medum masala chai"
Set-TextValue 16 4 "0322-7287568"
Set-TextValue 16 5 "Lahore"
Set-TextValue 16 6 "0322-7287568"
Set-TextValue 16 7 "cscdcd"
Set-TextValue 16 8 "csssdccs"
Set-TextValue 16 9 "cdcdscd"

$excel.CutCopyMode = 0
